$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 79
$newRow = 80

# Seed formatting for the new row by copying the previous row's formats for
# the two styled columns (A: bold/bordered index column, E: date/time
# column), then overwrite with the real values below.
$ws.Range("A" + $srcRow).Copy() | Out-Null
$ws.Range("A" + $newRow).PasteSpecial(-4122) | Out-Null
$ws.Range("E" + $srcRow).Copy() | Out-Null
$ws.Range("E" + $newRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 79
$ws.Cells.Item($newRow, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item($newRow, 3).Value = "premijer-liga-bih"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45242.54166666666
$ws.Cells.Item($newRow, 6).Value = "Zvijezda 09"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "GOSK Gabela"
$ws.Cells.Item($newRow, 9).Value = 2
$ws.Cells.Item($newRow, 10).Value = 2.07
$ws.Cells.Item($newRow, 11).Value = "12/11/2023 05:12"
$ws.Cells.Item($newRow, 12).Value = 2.05
$ws.Cells.Item($newRow, 13).Value = "12/11/2023 12:51"
$ws.Cells.Item($newRow, 14).Value = 3.38
$ws.Cells.Item($newRow, 15).Value = "12/11/2023 05:12"
$ws.Cells.Item($newRow, 16).Value = 3.07
$ws.Cells.Item($newRow, 17).Value = "12/11/2023 12:51"
$ws.Cells.Item($newRow, 18).Value = 3.23
$ws.Cells.Item($newRow, 19).Value = "12/11/2023 05:12"
$ws.Cells.Item($newRow, 20).Value = 3.87
$ws.Cells.Item($newRow, 21).Value = "12/11/2023 12:51"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/zvijezda-09-nk-gosk-gabela/tbYBHBhH/"
